# Update "想去人数" (want-to-go count) figures that changed between the
# previous data scrape and the latest one (gh-pages output regenerated at
# commit 456a3b4). The same underlying event data is duplicated across the
# "展览" and "全部类型" worksheets, so both need to be updated identically.

$wb = $excel.ActiveWorkbook

# row -> new F-column value, expressed relative to the "展览" sheet layout
$updates = @{
    2  = 7143
    3  = 20
    4  = 463
    7  = 173
    8  = 124
    11 = 57
    12 = 210
    13 = 9
    14 = 455
    16 = 1846
    17 = 47
    18 = 3
    19 = 3711
    23 = 34
    24 = 30
    25 = 2363
    26 = 18
    27 = 283
    31 = 7
    34 = 22
    36 = 1384
    37 = 126
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Range("F$row").Value = $updates[$row]
}

# The "全部类型" sheet has one extra row (row 6) before the others, so every
# row from 7 onward in "展览" maps to row+1 in "全部类型".
$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates.Keys) {
    if ($row -lt 6) {
        $targetRow = $row
    } else {
        $targetRow = $row + 1
    }
    $ws4.Range("F$targetRow").Value = $updates[$row]
}
